# The document contains three circular "message bubble" shapes
# (Oval 20, Oval 21, Oval 22) that currently have a solid blue
# outline/stroke drawn around the gradient-filled circle. This change
# removes that outline from all three ovals (the gradient fill itself is
# left untouched) - i.e. it turns each shape's line format off.
#
# Document.Shapes is indexed in the order the drawings appear in the
# file, and the three ovals are the first three shapes in this
# document, so Shapes.Item(1), Item(2) and Item(3) are exactly the three
# shapes this change touches.

$d = $word.ActiveDocument

$d.Shapes.Item(1).Line.Visible = $false
$d.Shapes.Item(2).Line.Visible = $false
$d.Shapes.Item(3).Line.Visible = $false
